# Commit: Sat, Jun 20, 2020  1:04:48 AM
#
# The table on slide 6 (the "SOURCES OF FINANCE" table, graphicFrame
# "Google Shape;127;p18") has its table style switched from the
# presentation's local custom style to a standard built-in PowerPoint
# table style:
#   {7E7FEFBD-0514-4F57-806C-0FF3B1448D2A}  ->  {2E4A10BF-E3A9-41EE-9957-481A13FC6722}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $tbl = $sh.Table
        $tbl.ApplyStyle("{2E4A10BF-E3A9-41EE-9957-481A13FC6722}")
        Write-Host "Applied table style {2E4A10BF-E3A9-41EE-9957-481A13FC6722} to shape $i ($($sh.Name)) on slide 6."
    }
}
